# The commit ("Fruta / hortaliza, semanal") refreshes the weekly price data:
# the data rows (2..33) of the single sheet are a straight re-shuffle of the
# same 32 records (prices/dates/quality bands swap places between rows; the
# set of values itself is unchanged). Row 15 keeps its original data.
#
# Strategy: snapshot every data row (columns A..R) with Value2 (plain COM
# property reads misbehave in this host for .Value), then write the rows
# back out in their new order from that snapshot so no row's data is lost
# while being overwritten mid-flight.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 33
$firstCol = 1
$lastCol = 18

# 1) Snapshot current contents of every data row.
$snapshot = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $rowVals = @{}
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $rowVals[$c] = $ws.Cells.Item($r, $c).Value2
    }
    $snapshot[$r] = $rowVals
}

# 2) Map: new row number -> source row number (from the snapshot above).
$rowMap = @{
    2  = 19
    3  = 20
    4  = 21
    5  = 17
    6  = 26
    7  = 8
    8  = 7
    9  = 16
    10 = 14
    11 = 32
    12 = 33
    13 = 28
    14 = 29
    15 = 15
    16 = 10
    17 = 11
    18 = 12
    19 = 13
    20 = 27
    21 = 5
    22 = 30
    23 = 31
    24 = 2
    25 = 3
    26 = 18
    27 = 4
    28 = 9
    29 = 6
    30 = 22
    31 = 23
    32 = 24
    33 = 25
}

# 3) Write each destination row from its mapped source row's snapshot.
foreach ($destRow in $rowMap.Keys) {
    $srcRow = $rowMap[$destRow]
    $srcVals = $snapshot[$srcRow]
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $ws.Cells.Item($destRow, $c).Value = $srcVals[$c]
    }
}
